$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '73.402.29'
$ws.Cells.Item(2, 5).Value = '  +0.08%  '

$ws.Cells.Item(3, 4).Value = '3.976.90'
$ws.Cells.Item(3, 5).Value = '  -2.05%  '

$ws.Cells.Item(4, 5).Value = '  -0.06%  '

$c = $ws.Cells.Item(5, 4)
$c.Value = "'608.39"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +6.88%  '

$c = $ws.Cells.Item(6, 4)
$c.Value = "'168.55"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +11.13%  '

$c = $ws.Cells.Item(7, 4)
$c.Value = "'0.682"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -2.06%  '

$ws.Cells.Item(8, 5).Value = '  -0.04%  '

$c = $ws.Cells.Item(9, 4)
$c.Value = "'0.786"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +2.05%  '

$ws.Cells.Item(10, 5).Value = '  +7.55%  '

$c = $ws.Cells.Item(11, 4)
$c.Value = "'56.01"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +3.92%  '

$c = $ws.Cells.Item(12, 4)
$c.Value = "'0.0000337"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +2.03%  '

$c = $ws.Cells.Item(13, 4)
$c.Value = "'11.29"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +1.65%  '

$ws.Cells.Item(14, 4).Value = '4.613.80'
$ws.Cells.Item(14, 5).Value = '  -2.09%  '

$ws.Cells.Item(15, 4).Value = '3.982.03'
$ws.Cells.Item(15, 5).Value = '  -2.32%  '

$c = $ws.Cells.Item(16, 4)
$c.Value = "'14.26"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -1.72%  '

$ws.Cells.Item(17, 5).Value = '  +1.75%  '

$c = $ws.Cells.Item(18, 4)
$c.Value = "'20.70"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.58%  '

$ws.Cells.Item(19, 4).Value = '73.354.04'
$ws.Cells.Item(19, 5).Value = '  +0.05%  '

$ws.Cells.Item(20, 5).Value = '  -1.36%  '

$c = $ws.Cells.Item(21, 4)
$c.Value = "'456.52"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.26%  '

$c = $ws.Cells.Item(22, 4)
$c.Value = "'4.82"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +6.75%  '

$c = $ws.Cells.Item(23, 4)
$c.Value = "'96.39"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.91%  '

$c = $ws.Cells.Item(24, 4)
$c.Value = "'3.43"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -3.79%  '

$c = $ws.Cells.Item(25, 4)
$c.Value = "'14.24"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -3.58%  '

$ws.Cells.Item(26, 5).Value = '  -1.38%  '

$c = $ws.Cells.Item(27, 4)
$c.Value = "'11.07"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -2.95%  '

$ws.Cells.Item(28, 5).Value = '  +0.26%  '

$c = $ws.Cells.Item(29, 4)
$c.Value = "'10.53"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -4.60%  '

$c = $ws.Cells.Item(30, 4)
$c.Value = "'36.43"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -2.19%  '

$c = $ws.Cells.Item(31, 4)
$c.Value = "'7.91"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.90%  '

$ws.Cells.Item(32, 5).Value = '  +2.20%  '

$c = $ws.Cells.Item(33, 4)
$c.Value = "'0.0000106"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +15.31%  '

$c = $ws.Cells.Item(34, 4)
$c.Value = "'0.130"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -3.74%  '

$c = $ws.Cells.Item(35, 4)
$c.Value = "'48.11"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -1.09%  '

$c = $ws.Cells.Item(36, 4)
$c.Value = "'70.83"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +3.96%  '

$c = $ws.Cells.Item(37, 4)
$c.Value = "'648.64"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -5.65%  '

$ws.Cells.Item(38, 5).Value = '  -3.83%  '

$ws.Cells.Item(39, 5).Value = '  +0.15%  '

$ws.Cells.Item(40, 2).Value = 'Kaspa'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Cells.Item(40, 4)
$c.Value = "'0.146"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -1.97%  '

$ws.Cells.Item(41, 2).Value = 'Dai'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Cells.Item(41, 4)
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.11%  '

$ws.Cells.Item(42, 5).Value = '  +0.00%  '

$ws.Cells.Item(43, 2).Value = 'dogwifhat'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Cells.Item(43, 4)
$c.Value = "'3.28"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +42.66%  '

$ws.Cells.Item(44, 2).Value = 'VeChain'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(44, 4)
$c.Value = "'0.0483"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -2.84%  '

$c = $ws.Cells.Item(45, 4)
$c.Value = "'10.59"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -5.57%  '

$ws.Cells.Item(46, 5).Value = '  -6.04%  '

$c = $ws.Cells.Item(47, 4)
$c.Value = "'0.149"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.42%  '

$c = $ws.Cells.Item(48, 4)
$c.Value = "'0.000302"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +8.02%  '

$c = $ws.Cells.Item(49, 4)
$c.Value = "'3.45"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +3.67%  '

$c = $ws.Cells.Item(50, 4)
$c.Value = "'2.58"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -5.85%  '

$ws.Cells.Item(51, 2).Value = 'Stacks'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Cells.Item(51, 4)
$c.Value = "'3.01"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -3.40%  '
